$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Resistor grouping / value re-assignment (rows 5-8)
#    R4 was reassigned into the 2.2K group (row 5), and the 1K / 4.7K / 10K
#    resistor rows were re-labelled: R3 stays 1K-group owner but now holds
#    1K directly (row6), R7/R8 move to 4.7K (row7, qty 2), R9 moves to 10K
#    (row8, qty 1).
# ---------------------------------------------------------------------------

# Row 5: 2.2K group now includes R4 instead of R7,R8,R9
$ws.Range("C5").Value = "R1,R2,R4,R5,R6"

# Row 6: now 1K (was 10K)
$ws.Range("D6").Value = "1K"
$ws.Range("E6").Value = "311-1.00KCRCT-ND"
$ws.Range("F6").Value = "RC0805FR-071KL"
$ws.Range("G6").Value = "RES 1K OHM 1% 1/8W 0805"

# Row 7: now 4.7K, quantity 2, designators R7,R8 (was 1K, qty 1, R4)
$ws.Range("B7").Value = 2
$ws.Range("C7").Value = "R7,R8"
$ws.Range("D7").Value = "4.7K"
$ws.Range("E7").Value = "311-4.7KARCT-ND"
$ws.Range("F7").Value = "RC0805JR-074K7L"
$ws.Range("G7").Value = "RES 4.7K OHM 5% 1/8W 0805"

# Row 8: now 10K, quantity 1, designator R9 (was 4.7K, qty 2, R5,R6)
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = "R9"
$ws.Range("D8").Value = "10K"
$ws.Range("E8").Value = "311-10.0KCRCT-ND"
$ws.Range("F8").Value = "RC0805FR-0710KL"
$ws.Range("G8").Value = "RES 10K OHM 1% 1/8W 0805"

# ---------------------------------------------------------------------------
# 2. Total price formula: operand order changed (value identical, commutative)
# ---------------------------------------------------------------------------
$ws.Range("I13").Formula = "=H2*B2+H3*B3+H4*B4+H5*B5+H8*B8+H6*B6+H7*B7+H9*B9+H10*B10"

# ---------------------------------------------------------------------------
# 3. Hyperlinks: the Digikey links must now point at the rows holding the
#    matching part after the re-shuffle above. The displayed text in J6:J8
#    must change to match the (moved) link target; J2:J5,J9:J10 keep their
#    existing text. Remove all existing hyperlinks and re-create them in the
#    correct final mapping, preserving relationship id order (rId1..rId9)
#    exactly as in the target workbook.
# ---------------------------------------------------------------------------
$ws.Range("A1").Hyperlinks.Delete()

$ws.Range("J7").Value = "https://www.digikey.ca/en/products/detail/yageo/RC0805JR-074K7L/728327"
$ws.Range("J6").Value = "https://www.digikey.ca/en/products/detail/yageo/RC0805FR-071KL/727444"
$ws.Range("J8").Value = "https://www.digikey.ca/en/products/detail/yageo/RC0805FR-0710KL/727535"

$ws.Hyperlinks.Add($ws.Range("J7"), "https://www.digikey.ca/en/products/detail/yageo/RC0805JR-074K7L/728327")
$ws.Hyperlinks.Add($ws.Range("J6"), "https://www.digikey.ca/en/products/detail/yageo/RC0805FR-071KL/727444")
$ws.Hyperlinks.Add($ws.Range("J8"), "https://www.digikey.ca/en/products/detail/yageo/RC0805FR-0710KL/727535")
$ws.Hyperlinks.Add($ws.Range("J9"), "https://www.digikey.ca/en/products/detail/avx-corporation/08053C104KAT2A/1116281")
$ws.Hyperlinks.Add($ws.Range("J10"), "https://www.digikey.ca/en/products/detail/lite-on-inc/LTST-C190KRKT/386817")
$ws.Hyperlinks.Add($ws.Range("J5"), "https://www.digikey.ca/en/products/detail/yageo/RC0805FR-072K2L/727676")
$ws.Hyperlinks.Add($ws.Range("J4"), "https://www.digikey.ca/en/products/detail/sparkfun-electronics/PRT-14417/7652746 ")
$ws.Hyperlinks.Add($ws.Range("J3"), "https://www.digikey.ca/en/products/detail/texas-instruments/TCA9534DWR/6566100 ")
$ws.Hyperlinks.Add($ws.Range("J2"), "https://www.digikey.ca/en/products/detail/cui-devices/SJ1-3513/738683 ")

# Adding hyperlinks stamps a fresh "hyperlink" cell style on each touched
# cell; restore the original style (shared by all J-column link cells) by
# copying formats from a cell that already carries it.
$ws.Range("J2").Copy()
$ws.Range("J2:J10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4. Restore the cursor/selection position as left by the author.
# ---------------------------------------------------------------------------
$ws.Range("E14").Select()
